{"js": "// \"estructura completa de informe\"\n// The last (empty) paragraph of the body -- the one styled \"Heading 2\"\n// (\"T\u00edtulo2\") sitting right before the final sectPr, with a manual\n// numbering override (ilvl 0 / numId 0) and a hanging indent -- loses all of\n// that paragraph formatting, becoming a bare empty paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length === 0) {\n  throw new Error(\"Document body has no paragraphs.\");\n}\n\n// Target the trailing empty \"Heading 2\" paragraph specifically (fall back to\n// the very last body paragraph if, for some reason, nothing matches).\nfor (const p of items) {\n  p.load(\"text,style\");\n}\nawait context.sync();\n\nlet target = null;\nfor (let i = items.length - 1; i >= 0; i--) {\n  const p = items[i];\n  if (p.text === \"\" && p.style === \"Heading 2\") {\n    target = p;\n    break;\n  }\n}\nif (!target) {\n  target = items[items.length - 1];\n}\n\n// Resetting the style to the document default clears the paragraph's\n// <w:pPr> entirely (style override, numbering override, and indentation\n// all go away), leaving a plain empty paragraph -- exactly the\n// <w:p>...</w:p> -> <w:p/> collapse in the diff.\ntarget.style = \"Normal\";\n\nawait context.sync();\n", "ps1": "# \"estructura completa de informe\"\n# The last (empty) paragraph of the document -- styled \"Heading 2\"\n# (\"T\u00edtulo2\"), with a manual numbering override (ilvl 0 / numId 0) and a\n# hanging indent -- loses all of that paragraph formatting, becoming a bare\n# empty paragraph, right before the final section break.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$target = $null\n\n# Walk back from the end looking for the empty \"Heading 2\" paragraph; fall\n# back to the very last paragraph in the document if nothing matches.\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -eq \"`r\" -and $p.Style.NameLocal -eq \"Heading 2\") {\n        $target = $p\n        break\n    }\n}\nif ($null -eq $target) {\n    $target = $d.Paragraphs.Last\n}\n\n# Resetting the style to the document default clears the paragraph's\n# formatting entirely (style override, numbering override, and indentation\n# all go away), leaving a plain empty paragraph.\n$target.Style = \"Normal\"\n"}
